$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 244, shifting the existing rows 244:276 down to 245:277.
$ws.Rows("244:244").Insert()

# Populate the newly inserted row 244 with the new weekly price record
# (same market/category constants as the surrounding rows, new date + stats).
$ws.Range("A244").Value = 6
$ws.Range("B244").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C244").Value = "Metropolitana"
$ws.Range("D244").Value = 44918
$ws.Range("E244").Value = 13
$ws.Range("F244").Value = 100112029
$ws.Range("G244").Value = "Orégano"
$ws.Range("H244").Value = "Sin especificar"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 43
$ws.Range("K244").Value = 16000
$ws.Range("L244").Value = 17000
$ws.Range("M244").Value = 16442
$ws.Range("N244").Value = "$/docena de atados"
$ws.Range("O244").Value = "Región Metropolitana"
$ws.Range("P244").Value = 5481
$ws.Range("Q244").Value = 3
$ws.Range("R244").Value = "Hortaliza"
